# Update NATMI LR-pair output (Anxa1-Dysf) with new TPM-derived expression values.
# Only the "ECs" cluster's ligand (Anxa1) and receptor (Dysf) average/total
# expression values changed upstream; all derived specificity and edge-weight
# columns are recomputed from those new base values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 2.811979666666667
$ws.Cells.Item(2, 8).Value2 = 8.435939000000001
$ws.Cells.Item(2, 9).Value2 = 0.01221198172659148
$ws.Cells.Item(2, 10).Value2 = 0.01221198172659148
$ws.Cells.Item(2, 13).Value2 = 34.417786
$ws.Cells.Item(2, 14).Value2 = 103.253358
$ws.Cells.Item(2, 15).Value2 = 0.8460109765801216
$ws.Cells.Item(2, 16).Value2 = 0.8460109765801216
$ws.Cells.Item(2, 17).Value2 = 96.78211440368469
$ws.Cells.Item(2, 18).Value2 = 871.039029633162
$ws.Cells.Item(2, 19).Value2 = 0.01033147058649226
$ws.Cells.Item(2, 20).Value2 = 0.01033147058649226
$ws.Cells.Item(3, 7).Value2 = 2.811979666666667
$ws.Cells.Item(3, 8).Value2 = 8.435939000000001
$ws.Cells.Item(3, 9).Value2 = 0.01221198172659148
$ws.Cells.Item(3, 10).Value2 = 0.01221198172659148
$ws.Cells.Item(3, 14).Value2 = 0.5243180000000001
$ws.Cells.Item(3, 15).Value2 = 0.004296022829771175
$ws.Cells.Item(3, 16).Value2 = 0.004296022829771176
$ws.Cells.Item(3, 17).Value2 = 0.4914571849557779
$ws.Cells.Item(3, 18).Value2 = 4.423114664602001
$ws.Cells.Item(3, 19).Value2 = 0.00005246295229418542
$ws.Cells.Item(3, 20).Value2 = 0.00005246295229418542
$ws.Cells.Item(4, 7).Value2 = 2.811979666666667
$ws.Cells.Item(4, 8).Value2 = 8.435939000000001
$ws.Cells.Item(4, 9).Value2 = 0.01221198172659148
$ws.Cells.Item(4, 10).Value2 = 0.01221198172659148
$ws.Cells.Item(4, 15).Value2 = 0.1496930005901073
$ws.Cells.Item(4, 16).Value2 = 0.1496930005901073
$ws.Cells.Item(4, 17).Value2 = 17.12460654719478
$ws.Cells.Item(4, 18).Value2 = 154.121458924753
$ws.Cells.Item(4, 19).Value2 = 0.001828048187805037
$ws.Cells.Item(4, 20).Value2 = 0.001828048187805038
$ws.Cells.Item(5, 9).Value2 = 0.8095640809678946
$ws.Cells.Item(5, 10).Value2 = 0.8095640809678947
$ws.Cells.Item(5, 13).Value2 = 34.417786
$ws.Cells.Item(5, 14).Value2 = 103.253358
$ws.Cells.Item(5, 15).Value2 = 0.8460109765801216
$ws.Cells.Item(5, 16).Value2 = 0.8460109765801216
$ws.Cells.Item(5, 17).Value2 = 6415.938482018797
$ws.Cells.Item(5, 18).Value2 = 57743.44633816917
$ws.Cells.Item(5, 19).Value2 = 0.6849000987438371
$ws.Cells.Item(5, 20).Value2 = 0.6849000987438372
$ws.Cells.Item(6, 9).Value2 = 0.8095640809678946
$ws.Cells.Item(6, 10).Value2 = 0.8095640809678947
$ws.Cells.Item(6, 14).Value2 = 0.5243180000000001
$ws.Cells.Item(6, 15).Value2 = 0.004296022829771175
$ws.Cells.Item(6, 16).Value2 = 0.004296022829771176
$ws.Cells.Item(6, 17).Value2 = 32.57997704069956
$ws.Cells.Item(6, 19).Value2 = 0.003477905774000795
$ws.Cells.Item(6, 20).Value2 = 0.003477905774000797
$ws.Cells.Item(7, 9).Value2 = 0.8095640809678946
$ws.Cells.Item(7, 10).Value2 = 0.8095640809678947
$ws.Cells.Item(7, 15).Value2 = 0.1496930005901073
$ws.Cells.Item(7, 16).Value2 = 0.1496930005901073
$ws.Cells.Item(7, 19).Value2 = 0.1211860764500567
$ws.Cells.Item(7, 20).Value2 = 0.1211860764500567
$ws.Cells.Item(8, 9).Value2 = 0.1782239373055139
$ws.Cells.Item(8, 10).Value2 = 0.1782239373055139
$ws.Cells.Item(8, 13).Value2 = 34.417786
$ws.Cells.Item(8, 14).Value2 = 103.253358
$ws.Cells.Item(8, 15).Value2 = 0.8460109765801216
$ws.Cells.Item(8, 16).Value2 = 0.8460109765801216
$ws.Cells.Item(8, 17).Value2 = 1412.456215211825
$ws.Cells.Item(8, 18).Value2 = 12712.10593690643
$ws.Cells.Item(8, 19).Value2 = 0.1507794072497922
$ws.Cells.Item(8, 20).Value2 = 0.1507794072497922
$ws.Cells.Item(9, 9).Value2 = 0.1782239373055139
$ws.Cells.Item(9, 10).Value2 = 0.1782239373055139
$ws.Cells.Item(9, 14).Value2 = 0.5243180000000001
$ws.Cells.Item(9, 15).Value2 = 0.004296022829771175
$ws.Cells.Item(9, 16).Value2 = 0.004296022829771176
$ws.Cells.Item(9, 17).Value2 = 7.172417751754224
$ws.Cells.Item(9, 18).Value2 = 64.551759765788
$ws.Cells.Item(9, 19).Value2 = 0.0007656541034761945
$ws.Cells.Item(9, 20).Value2 = 0.0007656541034761946
$ws.Cells.Item(10, 9).Value2 = 0.1782239373055139
$ws.Cells.Item(10, 10).Value2 = 0.1782239373055139
$ws.Cells.Item(10, 15).Value2 = 0.1496930005901073
$ws.Cells.Item(10, 16).Value2 = 0.1496930005901073
$ws.Cells.Item(10, 19).Value2 = 0.02667887595224553
$ws.Cells.Item(10, 20).Value2 = 0.02667887595224554
